$d = $word.ActiveDocument

# 1. Fix split run "análi" + "sis de riesgos..." -> merge into one run (no GoBack bookmark here anymore)
$d.Content.Find.Execute("el análisis de riesgos, y los planes de contingencia pensados por si ocurren estos riesgos. ", $true, $false, $false, $false, $false, $true, 1, $false, "el análisis de riesgos, y los planes de contingencia pensados por si ocurren estos riesgos. ", 2) | Out-Null

# 2. Spelling fixes
$d.Content.Find.Execute("librerias", $true, $false, $false, $false, $false, $true, 1, $false, "librerías", 2) | Out-Null
$d.Content.Find.Execute("erronea", $true, $false, $false, $false, $false, $true, 1, $false, "errónea", 2) | Out-Null
$d.Content.Find.Execute("requetimientos", $true, $false, $false, $false, $false, $true, 1, $false, "requerimientos", 2) | Out-Null
